# Update order list with new promoter targeting guides.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new Sequence value (column B)
$updates = @{
    4  = "GTGTCCTTACGGGTGCGTCC"  # ACP2_3
    10 = "TTTGCAGCCAGCCCCAAAGG"  # CLU_3
    12 = "GCCCGCAGAGCTGCCCTGAG"  # FERMT2_2
    14 = "CAGCTCCAGCTTCCCCACCC"  # RPA1|SMYD4_1
    15 = "GGACCATGGGTGGGTCACGT"  # RPA1|SMYD4_2
    16 = "ATGGTCGGCCAACTGAGCGA"  # RPA1|SMYD4_3
    19 = "AGCGCCTAAGCCCCGCCCCT"  # SNX1_3
    22 = "GCGCTCTACTCACCCGCGCG"  # TSPAN14_3
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
